$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (new random test-credential data)
$ws.Range("A2").Value = "QZOgc681"
$ws.Range("B2").Value = 2012454478
$ws.Range("C2").Value = "oglzrte35"
$ws.Range("D2").Value = "B2`$x!j7R"
$ws.Range("E2").Value = "MR"
$ws.Range("F2").Value = "DkQdKvAK"
$ws.Range("G2").Value = "oaSk"
$ws.Range("H2").Value = "Candidate"

# Add a brand new row 3 with a new set of test credentials
$ws.Range("A3").Value = "obSus287"
$ws.Range("B3").Value = 2012454479
$ws.Range("C3").Value = "ctxlhsa87"
$ws.Range("D3").Value = "u#4%SdC3"
$ws.Range("E3").Value = "MR"
$ws.Range("F3").Value = "FNLyCDag"
$ws.Range("G3").Value = "OmID"
$ws.Range("H3").Value = "Candidate"

# Match the row 2 style (borders/font) for the new row 3 cells
$ws.Range("A2:H2").Copy()
$ws.Range("A3:H3").PasteSpecial(-4122) | Out-Null

# Re-apply values since paste-special of formats shouldn't touch them,
# but ensure they remain correct
$ws.Range("A3").Value = "obSus287"
$ws.Range("B3").Value = 2012454479
$ws.Range("C3").Value = "ctxlhsa87"
$ws.Range("D3").Value = "u#4%SdC3"
$ws.Range("E3").Value = "MR"
$ws.Range("F3").Value = "FNLyCDag"
$ws.Range("G3").Value = "OmID"
$ws.Range("H3").Value = "Candidate"

$excel.CutCopyMode = 0
